# Import devices from excel 2
# Convert the "price" column (F) sample values from text placeholders
# (a5, b5, c5, e5) to actual numeric values, and update the active
# selection to F7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 10
$ws.Range("F3").Value = 20
$ws.Range("F4").Value = 30
$ws.Range("F6").Value = 40

$ws.Range("F7").Select()
